# Auto-generated edit script: updates Kraken_Profits leve-profit values
# across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets, per the scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 350
$ws.Range("I2").Value = 290
$ws.Range("J2").Value = 450
$ws.Range("K2").Value = 290
$ws.Range("L2").Value = 450
$ws.Range("M2").Value = -177
$ws.Range("N2").Value = -676

$ws.Range("H9").Value = 273.6154
$ws.Range("I9").Value = 150.81818
$ws.Range("K9").Value = 150.81818
$ws.Range("M9").Value = 18.18181999999999

$ws.Range("H29").Value = 400

$ws.Range("H38").Value = 698.9
$ws.Range("I38").Value = 9.75
$ws.Range("J38").Value = 1158.3334
$ws.Range("K38").Value = 29.25
$ws.Range("L38").Value = 3475.0002
$ws.Range("M38").Value = 342.75
$ws.Range("N38").Value = -4219.0002

$ws.Range("H41").Value = 1678.1111
$ws.Range("I41").Value = 2220
$ws.Range("J41").Value = 1000.75
$ws.Range("K41").Value = 2220
$ws.Range("L41").Value = 1000.75
$ws.Range("M41").Value = -1780
$ws.Range("N41").Value = -1880.75

$ws.Range("H43").Value = 3284
$ws.Range("J43").Value = 3447.75
$ws.Range("L43").Value = 3447.75
$ws.Range("N43").Value = -3585.75

$ws.Range("H69").Value = 7015
$ws.Range("J69").Value = 7015
$ws.Range("L69").Value = 21045
$ws.Range("N69").Value = -22793

$ws.Range("H72").Value = 7015
$ws.Range("J72").Value = 7015
$ws.Range("L72").Value = 63135
$ws.Range("N72").Value = -71871

$ws.Range("H107").Value = 2702
$ws.Range("I107").Value = 1100
$ws.Range("J107").Value = 5906
$ws.Range("K107").Value = 1100
$ws.Range("L107").Value = 5906
$ws.Range("M107").Value = 820
$ws.Range("N107").Value = -9746

$ws.Range("H132").Value = 3497.4783
$ws.Range("I132").Value = 2870
$ws.Range("K132").Value = 8610
$ws.Range("M132").Value = -6080

$ws.Range("H137").Value = 2995.182
$ws.Range("J137").Value = 3387.25
$ws.Range("L137").Value = 10161.75
$ws.Range("N137").Value = -15261.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2210
$ws.Range("I74").Value = 762.5
$ws.Range("K74").Value = 762.5
$ws.Range("M74").Value = 111.5

$ws.Range("H77").Value = 2210
$ws.Range("I77").Value = 762.5
$ws.Range("K77").Value = 3812.5
$ws.Range("M77").Value = 555.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4332.6665
$ws.Range("I20").Value = 4332.6665
$ws.Range("K20").Value = 4332.6665
$ws.Range("M20").Value = -4085.6665

$ws.Range("H64").Value = 547
$ws.Range("I64").Value = 399
$ws.Range("J64").Value = 695
$ws.Range("K64").Value = 399
$ws.Range("L64").Value = 695
$ws.Range("M64").Value = -174
$ws.Range("N64").Value = -1145

$ws.Range("H67").Value = 547
$ws.Range("I67").Value = 399
$ws.Range("J67").Value = 695
$ws.Range("K67").Value = 399
$ws.Range("L67").Value = 695
$ws.Range("M67").Value = 381
$ws.Range("N67").Value = -2255

$ws.Range("H95").Value = 22312
$ws.Range("J95").Value = 22312
$ws.Range("L95").Value = 22312
$ws.Range("N95").Value = -27804

$ws.Range("H107").Value = 1799
$ws.Range("I107").Value = 1835.75
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 1835.75
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = 84.25
$ws.Range("N107").Value = -5590

$ws.Range("H134").Value = 4413.1665
$ws.Range("I134").Value = 3196.25
$ws.Range("J134").Value = 6847
$ws.Range("K134").Value = 9588.75
$ws.Range("L134").Value = 20541
$ws.Range("M134").Value = -7053.75
$ws.Range("N134").Value = -25611

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 5001250
$ws.Range("I13").Value = 10000000
$ws.Range("J13").Value = 2500
$ws.Range("K13").Value = 10000000
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = -9999861
$ws.Range("N13").Value = -2778

$ws.Range("H31").Value = 4467.45
$ws.Range("I31").Value = 2327.3333
$ws.Range("J31").Value = 4845.1177
$ws.Range("K31").Value = 2327.3333
$ws.Range("L31").Value = 4845.1177
$ws.Range("M31").Value = -2032.3333
$ws.Range("N31").Value = -5435.1177

$ws.Range("H34").Value = 4467.45
$ws.Range("I34").Value = 2327.3333
$ws.Range("J34").Value = 4845.1177
$ws.Range("K34").Value = 2327.3333
$ws.Range("L34").Value = 4845.1177
$ws.Range("M34").Value = -2125.3333
$ws.Range("N34").Value = -5249.1177

$ws.Range("H96").Value = 4152.75
$ws.Range("J96").Value = 4152.75
$ws.Range("L96").Value = 4152.75
$ws.Range("N96").Value = -9644.75

$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H70").Value = 1454.5
$ws.Range("I70").Value = 1454.5
$ws.Range("K70").Value = 1454.5
$ws.Range("M70").Value = -1184.5

$ws.Range("H73").Value = 1454.5
$ws.Range("I73").Value = 1454.5
$ws.Range("K73").Value = 1454.5
$ws.Range("M73").Value = -518.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2214.2856
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -1705
$ws.Range("N22").Value = -3590

$ws.Range("H27").Value = 2214.2856
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -1893
$ws.Range("N27").Value = -3214

$ws.Range("H55").Value = 3314.7693
$ws.Range("I55").Value = 2949.6667
$ws.Range("J55").Value = 3627.7144
$ws.Range("K55").Value = 2949.6667
$ws.Range("L55").Value = 3627.7144
$ws.Range("M55").Value = -2776.6667
$ws.Range("N55").Value = -3973.7144

$ws.Range("H56").Value = 14025
$ws.Range("I56").Value = 14025
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 14025
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -13334
$ws.Range("N56").ClearContents()

$ws.Range("H132").Value = 10333
$ws.Range("I132").Value = 8000
$ws.Range("K132").Value = 24000
$ws.Range("M132").Value = -21470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 30085
$ws.Range("I58").Value = 30085
$ws.Range("K58").Value = 30085
$ws.Range("M58").Value = -29777

$ws.Range("H62").Value = 3499.3333
$ws.Range("I62").Value = 3249
$ws.Range("K62").Value = 3249
$ws.Range("M62").Value = -2625

$ws.Range("H65").Value = 3499.3333
$ws.Range("I65").Value = 3249
$ws.Range("K65").Value = 16245
$ws.Range("M65").Value = -13125

$ws.Range("H81").Value = 16105.9
$ws.Range("I81").Value = 25843.166
$ws.Range("K81").Value = 51686.332
$ws.Range("M81").Value = -50625.332

$ws.Range("H84").Value = 16105.9
$ws.Range("I84").Value = 25843.166
$ws.Range("K84").Value = 258431.66
$ws.Range("M84").Value = -253127.66

